# edit.ps1
# Applies two changes to the "<<TEL_EXT>>" / "Teléfono donde labora y/o
# extensión" table (the 9th table in the document, 1-based) per the diff:
#   1. Insert two empty centered paragraphs (Garet 10pt / es-MX rPr) right
#      before the paragraph holding "<<TEL_EXT>>" inside that cell.
#   2. Shrink the row height of the "Teléfono donde labora y/o extensión"
#      row from 602 twips (30.1pt) to 58 twips (2.9pt).

$d = $word.ActiveDocument
$t = $d.Tables.Item(9)

# --- Step 1: insert the two empty paragraphs ------------------------------
# Plain Range/Font property writes on a table cell's first paragraph were
# observed to bleed formatting into sibling paragraphs in this host, so the
# insertion + exact formatting is done atomically via a single InsertXML
# (flat-OPC) call that replaces the target paragraph with itself preceded
# by the two new empty paragraphs, all sharing the same pPr/rPr.
$cell = $t.Cell(1, 1)
$targetPara = $cell.Range.Paragraphs.Item(1)
$targetRange = $targetPara.Range

$xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Garet" w:eastAsia="Times New Roman" w:hAnsi="Garet" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-MX"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Garet" w:eastAsia="Times New Roman" w:hAnsi="Garet" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-MX"/></w:rPr></w:pPr></w:p><w:p w14:paraId="6C132D95" w14:textId="068A7978" w:rsidR="001D22C3" w:rsidRPr="001D22C3" w:rsidRDefault="00D42AFF" w:rsidP="001D22C3"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:eastAsia="es-MX"/></w:rPr></w:pPr><w:r w:rsidRPr="00E470EB"><w:rPr><w:rFonts w:ascii="Garet" w:eastAsia="Times New Roman" w:hAnsi="Garet" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-MX"/></w:rPr><w:t>&lt;&lt;TEL_EXT&gt;&gt;</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xmlPayload)

# --- Step 2: shrink row 2's height -----------------------------------------
$row = $t.Rows.Item(2)
$row.Height = 2.9
